$wb = $excel.ActiveWorkbook
Write-Host ($wb | Get-Member | Where-Object { $_.Name -like "*erson*" -or $_.Name -like "*omment*" })
